# Update the answers in the two-digit / one-digit division practice table.
# The table has 20 rows x 5 columns; only every 4th row (1, 5, 9, 13, 17 in
# 1-based indexing) actually contains an answer in each of its 5 cells.
# We address cells directly by (row, column) so that the update cannot be
# confused by duplicate text appearing elsewhere in the document (e.g. the
# new value for one cell equals the old value of another cell further down).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; New="95÷7=13, 4"},
    @{Row=1;  Col=2; New="36÷8=4, 4"},
    @{Row=1;  Col=3; New="51÷9=5, 6"},
    @{Row=1;  Col=4; New="58÷2=29, 0"},
    @{Row=1;  Col=5; New="53÷9=5, 8"},

    @{Row=5;  Col=1; New="53÷6=8, 5"},
    @{Row=5;  Col=2; New="50÷2=25, 0"},
    @{Row=5;  Col=3; New="13÷9=1, 4"},
    @{Row=5;  Col=4; New="20÷6=3, 2"},
    @{Row=5;  Col=5; New="55÷5=11, 0"},

    @{Row=9;  Col=1; New="57÷8=7, 1"},
    @{Row=9;  Col=2; New="92÷5=18, 2"},
    @{Row=9;  Col=3; New="86÷3=28, 2"},
    @{Row=9;  Col=4; New="78÷4=19, 2"},
    @{Row=9;  Col=5; New="23÷2=11, 1"},

    @{Row=13; Col=1; New="82÷7=11, 5"},
    @{Row=13; Col=2; New="23÷6=3, 5"},
    @{Row=13; Col=3; New="43÷2=21, 1"},
    @{Row=13; Col=4; New="77÷2=38, 1"},
    @{Row=13; Col=5; New="95÷3=31, 2"},

    @{Row=17; Col=1; New="24÷9=2, 6"},
    @{Row=17; Col=2; New="33÷5=6, 3"},
    @{Row=17; Col=3; New="39÷6=6, 3"},
    @{Row=17; Col=4; New="96÷3=32, 0"},
    @{Row=17; Col=5; New="53÷9=5, 8"}
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.New
}
